# Generate Report for Handoff
# Inserts a new "717fb483-3c34-4e0f-9a26-de255418033b.md" file entry (ready for
# handoff) ahead of the existing ".localization-config" row on every sheet
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Rows.Item(3).Insert()

$ov.Range("A3").Value = "717fb483-3c34-4e0f-9a26-de255418033b.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $null, $null, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $null, $null, "717fb483-3c34-4e0f-9a26-de255418033b.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Rows.Item(3).Insert()

$zh.Range("A3").Value = "717fb483-3c34-4e0f-9a26-de255418033b.md"
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-08 06:18:29"
$zh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $null, $null, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f80da01613880be026cb58b4d8a6dd21f9762c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.zh-cn.xlf", $null, $null, "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $null, $null, "717fb483-3c34-4e0f-9a26-de255418033b.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f80da01613880be026cb58b4d8a6dd21f9762c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf", $null, $null, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows.Item(3).Insert()

$de.Range("A3").Value = "717fb483-3c34-4e0f-9a26-de255418033b.md"
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf"
$de.Range("D3").Value = "2016-03-08 06:18:31"
$de.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/064e0f25-b0eb-4c64-af88-7a51c01369bc.md", $null, $null, "064e0f25-b0eb-4c64-af88-7a51c01369bc.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d8955e546467ebaa36cb5c84ac7b35b66c6f7ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.de-de.xlf", $null, $null, "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/e2e/717fb483-3c34-4e0f-9a26-de255418033b.md", $null, $null, "717fb483-3c34-4e0f-9a26-de255418033b.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d8955e546467ebaa36cb5c84ac7b35b66c6f7ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf", $null, $null, "717fb483-3c34-4e0f-9a26-de255418033b.e03fc1afdfa3174d27b56e0ffb00348886d138f6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/afe2a8caf478425c307531fd40f3dc0227bd18d1/.localization-config", $null, $null, ".localization-config")

$ov.Select()
$ov.Range("A1").Select()
